$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 ("8800001") was missing its "name" (col B) value during import, and
# had a spurious, leftover style on col D. Record the successful import:
# fill in the product name and drop the stray D8 formatting entirely.
$ws.Range("B8").Value = "Coisa Valida"
$ws.Range("D8").Clear()

# Move the selection/active cell down to B9 (view no longer needs to be
# scrolled to C1).
$ws.Range("B9").Select()
